$wb = $excel.ActiveWorkbook
$wsScenarios = $wb.Worksheets.Item("Scenarios")
$wsImprovements = $wb.Worksheets.Item("Improvements")

# --- "Scenarios" sheet (sheet1): rename a few scenario labels, update ---
# --- the "fuel for energy" row values ---
$wsScenarios.Range("A2").Value = "population growth"
$wsScenarios.Range("A3").Value = "change in electricity consumption per capita"
$wsScenarios.Range("A4").Value = "Change in desalinated water"
$wsScenarios.Range("A10").Value = "fuel for energy"

$wsScenarios.Range("D10").Value = 0
$wsScenarios.Range("H10").Value = 0.5
$wsScenarios.Range("I10").Value = 0.9

# --- "Improvements" sheet (sheet2): PV / area improvement value update ---
$wsImprovements.Range("B3").Value = 0.5

# --- Selection / active-sheet bookkeeping to match the authored session ---
$wsScenarios.Range("E10").Select()
$wsImprovements.Activate()
$wsImprovements.Range("D12").Select()
